# Apply the edits to the "VillaTest" sheet:
#  - B2 text "Shubham Vilas" -> "Shubham villas"
#  - B3 text remains "Savoy Suites Hotel Apartment" (unchanged)
#  - Active selection moves to B3

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VillaTest")
$ws.Activate()

$ws.Range("B2").Value = "Shubham villas"

$ws.Range("B3").Select()
